# Update "想去人数" (F column) counts to match the refreshed data snapshot
# (commit: "Update gh-pages to output generated at 456a3b4").
$wb = $excel.ActiveWorkbook

# 展览 (sheet1, 18 change(s))
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 774   # F2: 771 -> 774
$ws.Cells.Item(3, 6).Value = 14352   # F3: 14347 -> 14352
$ws.Cells.Item(8, 6).Value = 991   # F8: 990 -> 991
$ws.Cells.Item(14, 6).Value = 450   # F14: 449 -> 450
$ws.Cells.Item(15, 6).Value = 2134   # F15: 2135 -> 2134
$ws.Cells.Item(18, 6).Value = 921   # F18: 920 -> 921
$ws.Cells.Item(22, 6).Value = 828   # F22: 827 -> 828
$ws.Cells.Item(23, 6).Value = 3381   # F23: 3380 -> 3381
$ws.Cells.Item(25, 6).Value = 316   # F25: 315 -> 316
$ws.Cells.Item(26, 6).Value = 2443   # F26: 2441 -> 2443
$ws.Cells.Item(30, 6).Value = 1815   # F30: 1814 -> 1815
$ws.Cells.Item(31, 6).Value = 1084   # F31: 1083 -> 1084
$ws.Cells.Item(35, 6).Value = 4953   # F35: 4951 -> 4953
$ws.Cells.Item(36, 6).Value = 4920   # F36: 4919 -> 4920
$ws.Cells.Item(41, 6).Value = 3312   # F41: 3311 -> 3312
$ws.Cells.Item(47, 6).Value = 4449   # F47: 4448 -> 4449
$ws.Cells.Item(48, 6).Value = 619   # F48: 618 -> 619
$ws.Cells.Item(49, 6).Value = 301   # F49: 300 -> 301

# 演出 (sheet2, 1 change(s))
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 126   # F4: 125 -> 126

# 本地生活 (sheet3, 3 change(s))
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 7663   # F2: 7662 -> 7663
$ws.Cells.Item(3, 6).Value = 252   # F3: 251 -> 252
$ws.Cells.Item(4, 6).Value = 864   # F4: 862 -> 864

# 全部类型 (sheet4, 17 change(s))
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 774   # F2: 771 -> 774
$ws.Cells.Item(3, 6).Value = 252   # F3: 251 -> 252
$ws.Cells.Item(4, 6).Value = 864   # F4: 862 -> 864
$ws.Cells.Item(6, 6).Value = 14352   # F6: 14347 -> 14352
$ws.Cells.Item(10, 6).Value = 991   # F10: 0 -> 991
$ws.Cells.Item(11, 6).Value = 126   # F11: 125 -> 126
$ws.Cells.Item(15, 6).Value = 450   # F15: 449 -> 450
$ws.Cells.Item(18, 6).Value = 921   # F18: 920 -> 921
$ws.Cells.Item(21, 6).Value = 3381   # F21: 3380 -> 3381
$ws.Cells.Item(22, 6).Value = 316   # F22: 315 -> 316
$ws.Cells.Item(25, 6).Value = 1815   # F25: 1814 -> 1815
$ws.Cells.Item(33, 6).Value = 4953   # F33: 4951 -> 4953
$ws.Cells.Item(34, 6).Value = 4920   # F34: 4919 -> 4920
$ws.Cells.Item(39, 6).Value = 3312   # F39: 3311 -> 3312
$ws.Cells.Item(45, 6).Value = 4449   # F45: 4448 -> 4449
$ws.Cells.Item(46, 6).Value = 619   # F46: 618 -> 619
$ws.Cells.Item(47, 6).Value = 301   # F47: 300 -> 301
